# Regenerate Handback status report timestamps.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 707e14f8... row.
$wsOverview.Range("G3").Value = "2016-09-01 02:54:11"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 707e14f8... row.
$wsZhCn.Range("H3").Value = "2016-09-01 02:54:01"
$wsZhCn.Range("K3").Value = "2016-09-01 02:54:32"

# de-de sheet: "Correspond Handback DateTime" for the 707e14f8... row.
# (its "Correspond Handoff Datetime" at H3 shares the same string as
# Overview!G3 and is updated together with it above.)
$wsDeDe.Range("K3").Value = "2016-09-01 02:54:39"
